$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "327×2="
$t.Cell(1,2).Range.Text = "470×9="
$t.Cell(1,3).Range.Text = "778×6="
$t.Cell(1,4).Range.Text = "440×6="
$t.Cell(1,5).Range.Text = "693×5="
$t.Cell(5,1).Range.Text = "147×8="
$t.Cell(5,2).Range.Text = "493×6="
$t.Cell(5,3).Range.Text = "836×4="
$t.Cell(5,4).Range.Text = "110×2="
$t.Cell(5,5).Range.Text = "693×2="
$t.Cell(10,1).Range.Text = "910×8="
$t.Cell(10,2).Range.Text = "119×8="
$t.Cell(10,3).Range.Text = "681×2="
$t.Cell(10,4).Range.Text = "717×4="
$t.Cell(10,5).Range.Text = "725×5="
$t.Cell(15,1).Range.Text = "794×9="
$t.Cell(15,2).Range.Text = "924×2="
$t.Cell(15,3).Range.Text = "213×3="
$t.Cell(15,4).Range.Text = "812×8="
$t.Cell(15,5).Range.Text = "948×4="
$t.Cell(20,1).Range.Text = "837×4="
$t.Cell(20,2).Range.Text = "171×9="
$t.Cell(20,3).Range.Text = "404×4="
$t.Cell(20,4).Range.Text = "207×8="
$t.Cell(20,5).Range.Text = "972×7="

Write-Host "done"
